$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "image" column to Table1 (extends A1:F33 -> A1:G33)
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# Header + first data row values
$ws.Cells.Item(1, 7).Value = "image"
$ws.Cells.Item(2, 7).Value = "spirits/image1.png"

# Column width (approximate Excel's computed best-fit width for the new column)
$ws.Columns.Item(7).ColumnWidth = 15.5

# Match the recorded selection after the edit
$ws.Range("I5").Select() | Out-Null
